# More data model tweaks
#
# - Switch the active tab from "Inspections" (index 2) to "Replacement" (index 1).
# - On the "Replacement" sheet: scroll/select so I1 is the active cell of a
#   I1:L4 selection (mirrors the G1 top-left scroll position in the saved view).
# - Clear the I1:L1 header cells completely (cell + formatting) and blank out
#   (but keep the styled, empty cells for) I2:L4's values.
# - Activating "Replacement" naturally drops tabSelected from "Inspections".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Replacement")

# Make it the active sheet (this flips workbook activeTab + tabSelected).
$ws.Activate()

# Remove the now-unused header cells (I1:L1) entirely - style and value both go.
$ws.Range("I1:L1").Clear()

# Blank out the data cells beneath them, keeping their existing cell style.
$ws.Range("I2:L4").ClearContents()

# Scroll the view so column G is left-most...
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1

# ...and select I1:L4, with I1 as the active cell.
$ws.Range("I1:L4").Select() | Out-Null
